$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row 13 (inherits formatting from row 12 above, same as Excel's
# default "insert row" behaviour) and populate it with the new model's data.
$ws.Rows(13).Insert()

$ws.Range("A13").Value = "Combined Gemma Model"
$ws.Range("B13").Value = 0.0064874064200000003
$ws.Range("C13").Value = 8.94

$ws.Range("D13").Formula = "= -LOG(B13)"
$ws.Range("E13").Formula = "=TAN(3.14159265  * (  (C13 / 9)  - 11/18   )  )"
$ws.Range("F13").Formula = "=E13 + D13"

# Match the author's final selection.
$null = $ws.Range("B10").Select()
